$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update sheet name (date progression)
$ws.Name = "Through 2022-02-26"

# Update the "February (through 02-25)" label text
$ws.Range("A3").Value = "February (through 02-26)"

# Update row 3 (February) values for columns C through I
$ws.Range("C3").Value = 32
$ws.Range("D3").Value = 55
$ws.Range("E3").Value = 50
$ws.Range("F3").Value = 29
$ws.Range("G3").Value = 65
$ws.Range("H3").Value = 111
$ws.Range("I3").Value = 133

# Update row 4 (Total) values for columns C through I
$ws.Range("C4").Value = 83
$ws.Range("D4").Value = 130
$ws.Range("E4").Value = 136
$ws.Range("F4").Value = 78
$ws.Range("G4").Value = 139
$ws.Range("H4").Value = 328
$ws.Range("I4").Value = 292
